$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price-column (D) values look numeric ("4.71", "67.806.61", etc.)
# but must stay as TEXT, matching the source data which stores them as
# inline/shared strings. Setting NumberFormat="@" immediately before the
# value assignment prevents Excel from re-interpreting them as numbers,
# and resetting the Style back to "Normal" immediately after removes the
# temporary text-format styling so no stray style index is left behind.
# This must be done cell-by-cell (not on multi-area/union ranges) because
# NumberFormat/Style assignment on a union range only reliably affects
# the first contiguous area in this runtime.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.806.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.64%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.396.60'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.68%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '550.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.16%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.502'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.29%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.156'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.57%  '

$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.326'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '67.682.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.95%  '

$ws.Range("E14").Value = '  -1.03%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '328.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.90%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -5.20%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.59%  '

$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("E21").Value = '  -4.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.95'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0₃0789'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.47%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.96'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.83%  '

$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '413.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.25%  '

$ws.Range("E30").Value = '  -2.51%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.38%  '

$ws.Range("E33").Value = '  -0.02%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.62%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.104'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.290'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.43%  '

$ws.Range("E38").Value = '  -2.20%  '

$ws.Range("E39").Value = '  -5.77%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.25'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '127.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.55%  '

$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.92'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.91%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0701'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.55%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.470'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.551'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0909'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.20%  '

$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("E48").Value = '  -7.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.26%  '

$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0200'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.35%  '

$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0424'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.25%  '
